$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item('展览')
$ws1.Range('B4').Value = '2024-03-23'
$ws1.Range('C4').Value = '北京·thebONE×Ilike动漫游戏嘉年华S4'
$ws1.Range('D4').Value = '小关路39号 北投购物公园'
$ws1.Range('E4').Value = '2024.03.23 10:00-03.24 17:00'
$ws1.Range('F4').Value = 5843
$ws1.Range('G4').Value = 70
$ws1.Range('H4').Value = 'https://show.bilibili.com/platform/detail.html?id=79601'
$ws1.Range('I4').Value = '//i2.hdslb.com/bfs/openplatform/202401/pSrsMI9z1705646196593.jpeg'
$ws1.Range('C5').Value = '北京·万游引力S6 知名声优 张思王之内场见面签售会'
$ws1.Range('D5').Value = '半截塔路53号首创郎园station西门 郎园station中央车站文化中心'
$ws1.Range('E5').Value = '2024.03.23 10:00-03.23 17:00'
$ws1.Range('F5').Value = 66
$ws1.Range('G5').Value = 138
$ws1.Range('H5').Value = 'https://show.bilibili.com/platform/detail.html?id=81856'
$ws1.Range('I5').Value = '//i2.hdslb.com/bfs/openplatform/202402/OjoZZB5o1708226039086.jpeg'
$ws1.Range('C6').Value = '北京·万游引力国潮动漫嘉年华s6'
$ws1.Range('E6').Value = '2024.03.23 10:00-03.24 17:00'
$ws1.Range('F6').Value = 2880
$ws1.Range('G6').Value = 75
$ws1.Range('H6').Value = 'https://show.bilibili.com/platform/detail.html?id=79322'
$ws1.Range('I6').Value = '//i2.hdslb.com/bfs/openplatform/202402/wqACkjUk1708236212668.jpeg'
$ws1.Range('C7').Value = '北京·排球少年ONLY'
$ws1.Range('D7').Value = '永外高庄138号  大红门会展中心'
$ws1.Range('E7').Value = '2024.03.23 10:00-03.23 17:00'
$ws1.Range('F7').Value = 1250
$ws1.Range('G7').Value = 60
$ws1.Range('H7').Value = 'https://show.bilibili.com/platform/detail.html?id=80510'
$ws1.Range('I7').Value = '//i1.hdslb.com/bfs/openplatform/202401/wNTz3awE1704441972575.jpeg'
$ws1.Range('C8').Value = '北京·西山动漫游戏嘉年华·次元漫境冬日派对'
$ws1.Range('D8').Value = '石景山路68号 北京首钢会展中心'
$ws1.Range('E8').Value = '2024.03.23 09:00-03.24 17:00'
$ws1.Range('F8').Value = 4593
$ws1.Range('G8').Value = '不可售'
$ws1.Range('H8').Value = 'https://show.bilibili.com/platform/detail.html?id=76891'
$ws1.Range('I8').Value = '//i0.hdslb.com/bfs/openplatform/202311/VFTEz3C11701046733452.jpeg'
$ws1.Range('F9').Value = 382
$ws1.Range('F10').Value = 426
$ws1.Range('F11').Value = 105
$ws1.Range('F13').Value = 661
$ws1.Range('F14').Value = 145
$ws1.Range('F15').Value = 4135
$ws1.Range('F16').Value = 4135
$ws1.Range('F18').Value = 78
$ws1.Range('F19').Value = 87
$ws1.Range('F21').Value = 187
$ws1.Range('F22').Value = 50
$ws1.Range('F23').Value = 6158
$ws1.Range('F24').Value = 6158
$ws1.Range('F26').Value = 86
$ws1.Range('F28').Value = 421
$ws1.Range('F29').Value = 193
$ws1.Range('F30').Value = 440
$ws1.Range('F31').Value = 4615
$ws1.Range('F32').Value = 1588
$ws1.Range('F34').Value = 1736
$ws1.Range('F35').Value = 5800
$ws1.Range('F36').Value = 93
$ws1.Range('F38').Value = 75
$ws1.Range('F39').Value = 59
$ws1.Range('F40').Value = 3857
$ws1.Range('F41').Value = 79
$ws1.Range('F44').Value = 2378
$ws1.Range('F45').Value = 19
$ws1.Range('F49').Value = 256
$ws1.Range('F50').Value = 670
$ws1.Range('F51').Value = 13

$ws2 = $wb.Worksheets.Item('演出')
$ws2.Range('F3').Value = 180
$ws2.Range('F4').Value = 22
$ws2.Range('F5').Value = 84
$ws2.Range('F6').Value = 25
$ws2.Range('F11').Value = 12

$ws3 = $wb.Worksheets.Item('本地生活')
$ws3.Range('F2').Value = 1395

$ws4 = $wb.Worksheets.Item('全部类型')
$ws4.Range('F2').Value = 1395
$ws4.Range('F4').Value = 5843
$ws4.Range('F6').Value = 2880
$ws4.Range('F7').Value = 1250
$ws4.Range('F8').Value = 426
$ws4.Range('F9').Value = 105
$ws4.Range('F11').Value = 180
$ws4.Range('F12').Value = 661
$ws4.Range('F13').Value = 145
$ws4.Range('F14').Value = 4135
$ws4.Range('F15').Value = 4135
$ws4.Range('F17').Value = 78
$ws4.Range('F18').Value = 87
$ws4.Range('F20').Value = 187
$ws4.Range('F21').Value = 50
$ws4.Range('F22').Value = 6158
$ws4.Range('F23').Value = 6158
$ws4.Range('F25').Value = 86
$ws4.Range('F26').Value = 421
$ws4.Range('F27').Value = 193
$ws4.Range('F28').Value = 440
$ws4.Range('F29').Value = 84
$ws4.Range('F30').Value = 4615
$ws4.Range('F31').Value = 1588
$ws4.Range('F32').Value = 25
$ws4.Range('F34').Value = 1736
$ws4.Range('F36').Value = 5800
$ws4.Range('F37').Value = 93
$ws4.Range('F39').Value = 3857
$ws4.Range('F44').Value = 2378
$ws4.Range('F45').Value = 19
$ws4.Range('F49').Value = 256
$ws4.Range('F51').Value = 12
